$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 959.9167
$ws.Range("I19").Value = 838.6667
$ws.Range("K19").Value = 838.6667
$ws.Range("M19").Value = -663.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2648.75
$ws.Range("I70").Value = 1233.6666
$ws.Range("J70").Value = 3120.4443
$ws.Range("K70").Value = 3700.9998
$ws.Range("L70").Value = 9361.332900000001
$ws.Range("M70").Value = -3430.9998
$ws.Range("N70").Value = -9901.332900000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2648.75
$ws.Range("I73").Value = 1233.6666
$ws.Range("J73").Value = 3120.4443
$ws.Range("K73").Value = 3700.9998
$ws.Range("L73").Value = 9361.332900000001
$ws.Range("M73").Value = -2764.9998
$ws.Range("N73").Value = -11233.3329

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 5073.8335
$ws.Range("I96").Value = 4323.25
$ws.Range("K96").Value = 12969.75
$ws.Range("M96").Value = -11596.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4447.8887
$ws.Range("I98").Value = 4004.5
$ws.Range("K98").Value = 4004.5
$ws.Range("M98").Value = -2506.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 4447.8887
$ws.Range("I122").Value = 4004.5
$ws.Range("K122").Value = 12013.5
$ws.Range("M122").Value = -9563.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5675.4585
$ws.Range("I132").Value = 4618.304
$ws.Range("K132").Value = 13854.912
$ws.Range("M132").Value = -11324.912

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 29681.81
$ws.Range("I135").Value = 971.1429000000001
$ws.Range("J135").Value = 119003.89
$ws.Range("K135").Value = 8740.286100000001
$ws.Range("L135").Value = 1071035.01
$ws.Range("M135").Value = -6205.286100000001
$ws.Range("N135").Value = -1076105.01

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9143.178
$ws.Range("I32").Value = 6147.1353
$ws.Range("K32").Value = 6147.1353
$ws.Range("M32").Value = -5860.1353

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 348196.44
$ws.Range("I45").Value = 910277.75
$ws.Range("J45").Value = 4702.3335
$ws.Range("K45").Value = 910277.75
$ws.Range("L45").Value = 4702.3335
$ws.Range("M45").Value = -909900.75
$ws.Range("N45").Value = -5456.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2654.7727
$ws.Range("I110").Value = 2606.2942
$ws.Range("J110").Value = 2819.6
$ws.Range("K110").Value = 2606.2942
$ws.Range("L110").Value = 2819.6
$ws.Range("M110").Value = -561.2941999999998
$ws.Range("N110").Value = -6909.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 22765
$ws.Range("I132").Value = 25977.191
$ws.Range("J132").Value = 5901
$ws.Range("K132").Value = 77931.573
$ws.Range("L132").Value = 17703
$ws.Range("M132").Value = -75401.573
$ws.Range("N132").Value = -22763

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1005.4
$ws.Range("I94").Value = 781.9231
$ws.Range("K94").Value = 781.9231
$ws.Range("M94").Value = -330.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2654.2354
$ws.Range("I99").Value = 1926.8334
$ws.Range("J99").Value = 4400
$ws.Range("K99").Value = 1926.8334
$ws.Range("L99").Value = 4400
$ws.Range("M99").Value = -428.8334
$ws.Range("N99").Value = -7396

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2813.2856
$ws.Range("I105").Value = 2389.1
$ws.Range("K105").Value = 2389.1
$ws.Range("M105").Value = -642.0999999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 174984.5
$ws.Range("J108").Value = 174984.5
$ws.Range("L108").Value = 174984.5
$ws.Range("N108").Value = -182664.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2411.8628
$ws.Range("I134").Value = 2223.7556
$ws.Range("J134").Value = 3822.6667
$ws.Range("K134").Value = 6671.266799999999
$ws.Range("L134").Value = 11468.0001
$ws.Range("M134").Value = -4136.266799999999
$ws.Range("N134").Value = -16538.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 27939.25
$ws.Range("J43").Value = 27939.25
$ws.Range("L43").Value = 27939.25
$ws.Range("N43").Value = -28307.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 27236.697
$ws.Range("I58").Value = 34078.516
$ws.Range("K58").Value = 34078.516
$ws.Range("M58").Value = -33875.516

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3845.4167
$ws.Range("J99").Value = 4530.875
$ws.Range("L99").Value = 4530.875
$ws.Range("N99").Value = -7526.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H101").Value = 27939.25
$ws.Range("J101").Value = 27939.25
$ws.Range("L101").Value = 27939.25
$ws.Range("N101").Value = -34429.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3845.4167
$ws.Range("J126").Value = 4530.875
$ws.Range("L126").Value = 13592.625
$ws.Range("N126").Value = -18532.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 27236.697
$ws.Range("I136").Value = 34078.516
$ws.Range("K136").Value = 102235.548
$ws.Range("M136").Value = -99685.54800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 923.6667
$ws.Range("I8").Value = 923.6667
$ws.Range("K8").Value = 2771.0001
$ws.Range("M8").Value = -2632.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1024.625
$ws.Range("I18").Value = 1024.625
$ws.Range("K18").Value = 3073.875
$ws.Range("M18").Value = -2904.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6683270.5
$ws.Range("I131").Value = 1818.8
$ws.Range("J131").Value = 10023996
$ws.Range("K131").Value = 5456.4
$ws.Range("L131").Value = 30071988
$ws.Range("M131").Value = -416.3999999999996
$ws.Range("N131").Value = -30082068

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5915.8335
$ws.Range("I70").Value = 4832.6665
$ws.Range("K70").Value = 4832.6665
$ws.Range("M70").Value = -4562.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5915.8335
$ws.Range("I73").Value = 4832.6665
$ws.Range("K73").Value = 4832.6665
$ws.Range("M73").Value = -3896.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 92443.09
$ws.Range("I107").Value = 143785.58
$ws.Range("K107").Value = 143785.58
$ws.Range("M107").Value = -141865.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4054.9524
$ws.Range("I122").Value = 2824.8572
$ws.Range("J122").Value = 4670
$ws.Range("K122").Value = 8474.571599999999
$ws.Range("L122").Value = 14010
$ws.Range("M122").Value = -6024.571599999999
$ws.Range("N122").Value = -18910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1024.2727
$ws.Range("I16").Value = 876.7
$ws.Range("K16").Value = 876.7
$ws.Range("M16").Value = -706.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5533.3335
$ws.Range("I40").Value = 5010
$ws.Range("K40").Value = 5010
$ws.Range("M40").Value = -4874

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1887.68
$ws.Range("I61").Value = 1903.875
$ws.Range("K61").Value = 1903.875
$ws.Range("M61").Value = -1701.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1887.68
$ws.Range("I113").Value = 1903.875
$ws.Range("K113").Value = 1903.875
$ws.Range("M113").Value = 266.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4519.231
$ws.Range("I122").Value = 3611.111
$ws.Range("K122").Value = 10833.333
$ws.Range("M122").Value = -8383.332999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1588.2084
$ws.Range("I122").Value = 1316.15
$ws.Range("K122").Value = 3948.45
$ws.Range("M122").Value = -1498.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 56008.668
$ws.Range("I132").Value = 54400.145
$ws.Range("K132").Value = 163200.435
$ws.Range("M132").Value = -160670.435

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4543.76
$ws.Range("I136").Value = 4884.1875
$ws.Range("K136").Value = 14652.5625
$ws.Range("M136").Value = -12102.5625
